# Apply "calibration changes/testing final data" edit to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New FVC (measured) values in column B ---
$bValues = @(2.9658000000000002, 2.9298999999999999, 2.9908000000000001, 3.0436999999999999, 2.9897, 3.004, 3.0196000000000001, 2.8353000000000002, 3.0272999999999999, 3.0329999999999999)

# --- New FEV (measured) values in column D ---
$dValues = @(2.5562999999999998, 2.1509, 2.0684, 2.2814000000000001, 1.7674000000000001, 1.9328000000000001, 2.1667000000000001, 2.5405000000000002, 2.1536, 2.0937000000000001)

# --- New calibrated FEV values in column H (was column G's formula target, now literal values) ---
$hValues = @(1.5348999999999999, 1.5204, 1.5109999999999999, 1.5168999999999999, 1.5091000000000001, 1.5362, 1.5208999999999999, 1.4976, 1.5338000000000001, 1.5072000000000001)

for ($i = 0; $i -lt 10; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
    $ws.Cells.Item($row, 4).Value = $dValues[$i]
}

# --- Move the "Trial" + "FEV (measured)" mini-table from columns F:H to G:H, clear old F/G/H ---
# Trial numbers are 1..10, same as column A
$trialNums = @(1, 2, 3, 4, 5, 6, 7, 8, 9, 10)

# Clear old F1:F11, G1, H1:H11 (the old layout)
$ws.Range("F1:H11").ClearContents()

# New header row: G1 = "Trial", H1 = "FEV (measured)", I1 = "FEV % Error"
$ws.Cells.Item(1, 7).Value = "Trial"
$ws.Cells.Item(1, 8).Value = "FEV (measured)"
$ws.Cells.Item(1, 9).Value = "FEV % Error"

# New data: G = trial number, H = literal calibrated value, I = % error formula
for ($i = 0; $i -lt 10; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $trialNums[$i]
    $ws.Cells.Item($row, 8).Value = $hValues[$i]
    $ws.Cells.Item($row, 9).Formula = "=ABS(1.5-H$row)/1.5*100"
}

# --- Column widths for new columns H (8) and I (9) ---
# (ColumnWidth is quantized internally; these inputs land on the closest
#  achievable grid point to the target stored widths 12.85546875 / 10.85546875)
$ws.Columns.Item(8).ColumnWidth = 12.022135416666666
$ws.Columns.Item(9).ColumnWidth = 10.022135416666666

# --- Update selection to D4 ---
$ws.Range("D4").Select()

$ws.Calculate()
